# BlackList.xlsx update
#
# The "Import" cell (B2) used to list two POJO imports
# ("com.redhat.prudential_poc.pojo.Application,com.redhat.prudential_poc.pojo.Insured")
# and is trimmed down to just the Insured import. The user's selection also
# moves from C9 to B3.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "com.redhat.prudential_poc.pojo.Insured"

$ws.Range("B3").Select() | Out-Null
